# Fixed Grogu information: update the "First Screen Appearance" text for
# Grogu (row 6) from "The Mandalorian (Season 1, 2019)" to
# "The Mandalorian: Season 1 (2019)".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sw_character_data")

$ws.Range("F6").Value = "The Mandalorian: Season 1 (2019)"
